# Update Logging_OKANT worksheet:
# - bump "Förändrad" (column C) date on all data rows from 45410 to 45412
# - delete the last row (row 29) which is no longer tracked
# - the row above the deleted one (row 28) loses its explicit custom row height

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 3)
    if ($c.Value2 -eq 45410) {
        $c.Value = 45412
    }
}

$ws.Rows.Item(29).Delete()
$ws.Rows.Item(28).AutoFit()
